$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new columns before column B, shifting existing B:V data to K:AE.
$ws.Range("B1:J1").EntireColumn.Insert()

# New weekly date headers for the newly inserted columns (most recent first).
$dates = @("Sep_08","Aug_25","Aug_04","Jul_23","Jul_17","Jul_07","Jun_30","Jun_24","Jun_16")
for ($i = 0; $i -lt $dates.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(1, $col).Value = $dates[$i]
}

# Fill the new columns with "UN" for every data row (rows with a label in column A).
$lastRow = 33
for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 1).Value()
    if ($label -ne $null -and $label -ne "") {
        for ($c = 2; $c -le 10; $c++) {
            $ws.Cells.Item($r, $c).Value = "UN"
        }
    }
}
